$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cea2.mex function now works without writing a new file; mark the
# corresponding SUBS rows (OUT1, NEWOF, ROCKET) as unused ("X") in column C,
# matching the other already-flagged rows.
$ws.Range("C16").Value = "X"
$ws.Range("C17").Value = "X"
$ws.Range("C19").Value = "X"

# Update the active selection left over from editing, as recorded by Excel.
$ws.Range("F16").Select()
